$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 1 table (rows 3-12) ---
$ws.Range("J3").Value = "Prazo Final:"

$ws.Range("B4").Value = "Estipular o tema e criar o Termo de Abertura do Projeto"
$ws.Range("B5").Value = "Configurar e Organizar o GitHub e Ferramenta de Gestão do Projeto`n"
$ws.Range("B7").Value = 'Desenvolver uma Camada de View Basica para a Pagina: "game_screen"'
$ws.Range("B8").Value = "Desenvolver as Regras de Negocio da Camada do Jogo: ""input_user_keyboard"", ""model_game""`ne ""render_screen"" "
$ws.Range("B9").Value = 'Relacionar os arquivos do diretorio "game_layer" com o controller e depois com o server.js'
$ws.Range("B10").Value = "Fazer Bateria de Testes Unitarios na Camada Jogo"

$ws.Range("B12").Value = "Finalizar toda Documentação inicial e a Camada do Jogo"
$ws.Range("D12").Value = "MAX"
$ws.Range("E12").Value = "MAX"

# --- Sprint 2 table (rows 15-24) ---
$ws.Range("B16").Value = 'Desenvolver a Camada de View das Paginas: "index",  "login_page" e "register_page"'
$ws.Range("B17").Value = 'Desenvolver Principais  as Regras de Negocio da camada Cadastro: "input_user_register" e "model_register"'
$ws.Range("B18").Value = 'Desenvolver as Regras de Negocio da camada Login: "input_user_login" e "model_login"'
$ws.Range("B19").Value = 'Relacionar os arquivos dos diretorios "user_register_layer" e "user_login_layer" com o controller e depois com o server.js'
$ws.Range("B20").Value = "Criar o banco de dados em seu modelo lógico, Dicionario de Dados e Banco Físico"
$ws.Range("B21").Value = "Fazer Bateria de Testes Unitarios nas Camadas de Cadastro e Login de Usuarios"
$ws.Range("B22").Value = "Fazer Bateria de Testes de Integração entre as 3 Camadas: Jogo, Login e Cadastro"
$ws.Range("B23").Value = "Fazer o Documento de Especificação Tecnica da Solução"

$ws.Range("B24").Value = "Finalizar toda Parte de Cadastro, Login, Banco de Dados e os Documentos Intermediarios "
$ws.Range("D24").Value = "MAX"
$ws.Range("E24").Value = "MAX"

# --- Sprint 3 table (rows 26-38) ---
$ws.Range("B26").Value = "Sprint 3 - Interligar os dados do Site com o BD"

$ws.Range("B28").Value = "Fazer o Model com o Mapeamento das Tabelas (entidades) em classes JavaScript "
$ws.Range("B29").Value = "Integração da API com Sistema "
$ws.Range("B30").Value = "Integração da aplicação web com o Banco de Dados"
$ws.Range("B31").Value = "Fazer Bateria de Testes Unitarios "
$ws.Range("B32").Value = "Fazer Bateria de Testes de Integração entre: Aplicação, API e Banco de Dados"
$ws.Range("B33").Value = 'Fazer Testes "END-TO-END" '
$ws.Range("B34").Value = "Criar uma POC utilizando Métricas Aplicadas aos Dados"
$ws.Range("B35").Value = 'Desenvolver um Fluxograma em linha de prioridade de "Atendimento do Suporte"  '
$ws.Range("B36").Value = "Criar o Arquivo PPT da Apresentação. Lembrando dos Temas pedidos na Apresentação pela Materia Socioemocional "

$ws.Range("A37").Value = "RF10"
$ws.Range("B37").Value = "Apresentar o Projeto para a Banca"

$ws.Range("A38").Value = "RF11"
$ws.Range("B38").Value = 'Concluir Projeto - "Space-Invaders-Bullet-Hell"'
$ws.Range("D38").Value = "MAX"
$ws.Range("E38").Value = "MAX"

# --- Row 4 formatting: pick up the same banded fill used by rows 5-11 ---
$ws.Range("B5:F5").Copy() | Out-Null
$ws.Range("B4:F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Reset the scroll position of the sheet view (was parked at A25) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
